$wb = $excel.ActiveWorkbook

$wsProps = $wb.Worksheets.Item("NP_Props")

# Convert "yes"/"no" text values in H2:I10 to numeric 1/0
for ($r = 2; $r -le 10; $r++) {
    $hVal = $wsProps.Cells.Item($r, 8).Value2
    if ($hVal -eq "yes") {
        $wsProps.Cells.Item($r, 8).Value = 1
    } elseif ($hVal -eq "no") {
        $wsProps.Cells.Item($r, 8).Value = 0
    }

    $iVal = $wsProps.Cells.Item($r, 9).Value2
    if ($iVal -eq "yes") {
        $wsProps.Cells.Item($r, 9).Value = 1
    } elseif ($iVal -eq "no") {
        $wsProps.Cells.Item($r, 9).Value = 0
    }
}

# Add new incubation concentration values for rows 3-6
$wsProps.Cells.Item(3, 10).Value = 4
$wsProps.Cells.Item(4, 10).Value = 4
$wsProps.Cells.Item(5, 10).Value = 4
$wsProps.Cells.Item(6, 10).Value = 4

# Make NP_Props the active sheet and set its selection
$wsProps.Activate()
$wsProps.Range("A2:A10").Select()
